$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (LoginName) on rows 2-6 gets a "1" suffix appended.
$ws.Range("C2").Value = "Pford1"
$ws.Range("C3").Value = "arthur1"
$ws.Range("C4").Value = "Zaphod1"
$ws.Range("C5").Value = "Gdirk1"
$ws.Range("C6").Value = "Trillian1"
